# [IMP] New data for test environment
# Update the sale-order test fixture: bump reference numbers / dates from the
# 2020/2021 test data to 2022, and refresh a couple of view-related cosmetics
# (tab ratio, page-layout zoom, active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Book-level view: widen the sheet-tab area a bit (tabRatio 500 -> 600 / 0.5 -> 0.6)
$excel.ActiveWindow.TabRatio = 0.6

# Refresh the sample order references / dates for the new test data
$ws.Range("D2").Value = "P1/2022/0001"
$ws.Range("D3").Value = 220123
$ws.Range("D4").Value = 22011214
$ws.Range("D5").Value = "IT/22/004"
$ws.Range("D8").Value = "P1/2022/0007"

# Sheet view: page-layout zoom 100 -> 60, and reset the active selection to A1
$excel.ActiveWindow.Zoom = 60
$ws.Range("A1").Select()
